# Project Journal.xlsx edit:
#   - Log a new entry (row 8) on the "Jesse" sheet for the work described in
#     the commit message.
#   - "Jesse" becomes the active tab (was "Main").
# Everything else (the "Main" sheet's totals, sharedStrings bookkeeping,
# sheet dimensions, etc.) is recomputed/maintained by Excel automatically
# once the underlying data changes.

$wb = $excel.ActiveWorkbook
$wsJesse = $wb.Worksheets.Item("Jesse")

# --- Jesse sheet: add a new log entry in row 8 ---
# A8: date (11/29/2017, serial 43068) -- same style/number-format as the
# rows above it (column A already carries style 6 = date format).
$wsJesse.Range("A8").Value = 43068

# B8: minutes spent. The "Total Time Spent (minutes)" cell (C2) on this
# sheet is `=SUM(B4:B200)`, so this also drives Main!B2 (`=(Jesse!C2)/60`)
# to recompute automatically once this value is written.
$wsJesse.Range("B8").Value = 90

# C8: description of the work done, taken from the commit message.
$wsJesse.Range("C8").Value = "Added data member 'shift' and functions 'setShift', 'getShift', and 'getNext' to Object superclass. `nEdited Object assignment in Rooms.h to account for shift data members.`nAdded shift values to room files for Doors and Windows. `nAdded description for 'shift' in BuildingInfo.txt.`nRemoved overloaded output operator for class Rooms. `nAdded function 'displayRoom' to class Rooms. "

# The new row's text wraps onto several lines, so its height is taller than
# the default row height.
$wsJesse.Rows.Item(8).RowHeight = 156.75

# --- View state: Jesse becomes the selected/active tab, scrolled so the new
# entry (row 7 onward) is visible ---
$wsJesse.Select()
$wsJesse.Range("G7").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
